# Fruta / hortaliza, semanal
# Insert a new weekly record as row 156 in the Piña - Vega Monumental
# Concepción price sheet, shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 156 (existing rows 156..217 shift to 157..218)
$ws.Rows("156:156").Insert()

# Populate the newly inserted row with this week's data
$ws.Range("A156").Value = 11
$ws.Range("B156").Value = 'Vega Monumental Concepción'
$ws.Range("C156").Value = 'Bíobío'
$ws.Range("D156").Value = 44876
$ws.Range("E156").Value = 8
$ws.Range("F156").Value = 'Fruta'
$ws.Range("G156").Value = 100108
$ws.Range("H156").Value = 'Tropicales y subtropicales'
$ws.Range("I156").Value = 100108005
$ws.Range("J156").Value = 'Piña'
$ws.Range("K156").Value = 'Caramelo'
$ws.Range("L156").Value = 'Segunda'
$ws.Range("M156").Value = 220
$ws.Range("N156").Value = 25000
$ws.Range("O156").Value = 27000
$ws.Range("P156").Value = 26091
$ws.Range("Q156").Value = '$/caja 14 unidades'
$ws.Range("R156").Value = 'Ecuador'
$ws.Range("S156").Value = 1864
$ws.Range("T156").Value = 14
